# Generate Report for Handback
# Updates the handoff/handback timestamps for the second file
# (4f38d11f-f14e-400f-b6ba-7c0ecfb544fc.md) across the Overview, zh-cn and
# de-de sheets, reflecting a freshly generated handback report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for row 3 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-15 14:46:20"

# --- zh-cn sheet: row 3 handoff / handback datetimes ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-15 14:46:15"
$wsZhCn.Range("K3").Value = "2016-08-15 14:46:31"

# --- de-de sheet: row 3 handoff / handback datetimes ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-15 14:46:20"
$wsDeDe.Range("K3").Value = "2016-08-15 14:46:39"
